# Remove the detail columns (Village/Town/City, Street, Line No, Pole No,
# Latitude, Longitude, FLOC) for pole records 50-117, leaving only the
# Item No (column A) in each of those rows. This mirrors the commit that
# strips the per-pole detail out of "page 3".."page 7" while keeping the
# numbering column intact.

$wb = $excel.ActiveWorkbook
$cols = @("B","D","G","H","I","J","K")

# Map of sheet name -> list of row numbers whose detail cells must be cleared.
$targets = @{
    "page 3" = @(34,36,38,40,42,44);
    "page 4" = @(6,8,10,12,14,16,18,20,22,24,26,28,30,32,34,36,38,40,42,44);
    "page 5" = @(6,8,10,12,14,16,18,20,22,24,26,28,30,32,34,36,38,40,42,44);
    "page 6" = @(6,8,10,12,14,16,18,20,22,24,26,28,30,32,34,36,38,40,42,44);
    "page 7" = @(6,8);
}

foreach ($sheetName in $targets.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $targets[$sheetName]) {
        foreach ($col in $cols) {
            $ws.Range("$col$r").Value = ""
        }
    }
}
